# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with the latest scraped figures, and swap the Maker / VeChain rows (49-50)
# to reflect their new relative ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (unstyled "Normal" cell) used to restore D-column cells to
# their original (unstyled) appearance after writing them as text - setting
# a numeric-looking string via .Value can otherwise flip the cell to a
# Number type / pick up a quote-prefix style; re-applying this style keeps
# both the underlying cell style index AND the stored type identical to a
# plain text write.
$normalStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "'58.948.70"
$ws.Range("D2").Style = $normalStyle
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "'2.627.09"
$ws.Range("D3").Style = $normalStyle
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'521.55"
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "'145.29"
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = $normalStyle
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "'2.640.22"
$ws.Range("D9").Style = $normalStyle
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "'3.087.69"
$ws.Range("D14").Style = $normalStyle
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "'58.909.97"
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "'20.88"
$ws.Range("D16").Style = $normalStyle
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "'0.0000136"
$ws.Range("D17").Style = $normalStyle
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'2.635.06"
$ws.Range("D18").Style = $normalStyle
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").Value = "'345.37"
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = $normalStyle
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'61.71"
$ws.Range("D24").Style = $normalStyle
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "'0.165"
$ws.Range("D26").Style = $normalStyle
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").Value = "'0.995"
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "'0.0₃0800"
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "'7.10"
$ws.Range("D29").Style = $normalStyle
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = $normalStyle
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'6.25"
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = "  +3.00%  "
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").Value = "'18.88"
$ws.Range("D33").Style = $normalStyle
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "'150.54"
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "'0.981"
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("D36").Value = "'3.98"
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").Value = "'1.42"
$ws.Range("D41").Style = $normalStyle
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").Value = "'278.42"
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = "  -5.25%  "
$ws.Range("D43").Value = "'0.995"
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").Value = "'0.607"
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "'19.44"
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("E47").Value = "  -3.35%  "
$ws.Range("D51").Value = "'4.62"
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = "  -1.79%  "

# Rows 49 and 50 swap order (Maker now ranks above VeChain), with refreshed
# Price / Volume(1h) figures for both coins.
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'1.990.71"
$ws.Range("D49").Style = $normalStyle
$ws.Range("E49").Value = "  +3.05%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0229"
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = "  -0.33%  "
